$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Activate()

# Delete entire row 17 ("Deletion" row), shifting row 18 (simulation_timepoints) up to row 17
$ws.Rows("17:17").Delete()

# Update the selection to match the post-edit state (A17:XFD17 selected, active cell A17)
$ws.Range("A17:XFD17").Select()
